$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42/43 swap: Maker <-> TrustWalletToken (with updated values)
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

# Update Price (D) and Volume(1h) (E) columns with new scraped values.
# Price values are stored as text in the workbook; cells whose new value would
# otherwise be auto-detected as a number are forced back to text format first
# so Excel keeps storing them as strings (matching the source data feed).
$ws.Range("D2").Value = '26.153.27'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.656.31'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.64'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5239'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2662'
$ws.Range("E8").Value = '  +1.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06355'
$ws.Range("E9").Value = '  +0.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.59'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07687'
$ws.Range("E11").Value = '  -1.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.610'
$ws.Range("E12").Value = '  +2.49%  '
$ws.Range("D13").Value = '1.677.48'
$ws.Range("E13").Value = '  +1.21%  '
$ws.Range("D14").Value = '1.884.21'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").Value = '0.0₅8205'
$ws.Range("E16").Value = '  +2.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.52'
$ws.Range("D18").Value = '26.144.76'
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.660'
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.56'
$ws.Range("E21").Value = '  +4.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '192.79'
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.955'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.56'
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.266'
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.515'
$ws.Range("E29").Value = '  +1.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05454'
$ws.Range("E30").Value = '  -4.28%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.466'
$ws.Range("E32").Value = '  -0.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.372'
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.565'
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9538'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.779'
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.402'
$ws.Range("E37").Value = '  -0.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5688'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01587'
$ws.Range("E39").Value = '  -0.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.875'
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8316'
$ws.Range("E42").Value = '  -1.67%  '
$ws.Range("D43").Value = '1.024.16'
$ws.Range("E43").Value = '  -3.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.33'
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("D45").Value = '1.795.17'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.81'
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("E47").Value = '  +6.01%  '
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.026'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("E50").Value = '  -1.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05203'
$ws.Range("E51").Value = '  -3.53%  '
